$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 4271.5
$ws.Range("J6").Value = 6833.3335
$ws.Range("L6").Value = 20500.0005
$ws.Range("N6").Value = -20724.0005
$ws.Range("H15").Value = 294.48648
$ws.Range("I15").Value = 294.48648
$ws.Range("K15").Value = 883.4594399999999
$ws.Range("M15").Value = -714.4594399999999
$ws.Range("H34").Value = 13098.667
$ws.Range("I34").Value = 13098.667
$ws.Range("K34").Value = 13098.667
$ws.Range("M34").Value = -12895.667
$ws.Range("H36").Value = 13098.667
$ws.Range("I36").Value = 13098.667
$ws.Range("K36").Value = 13098.667
$ws.Range("M36").Value = -12383.667
$ws.Range("H62").Value = 41667836
$ws.Range("I62").Value = 62500504
$ws.Range("K62").Value = 62500504
$ws.Range("M62").Value = -62499880
$ws.Range("H65").Value = 41667836
$ws.Range("I65").Value = 62500504
$ws.Range("K65").Value = 312502520
$ws.Range("M65").Value = -312499400
$ws.Range("H70").Value = 146200.58
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 203880.8
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 611642.3999999999
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -612182.3999999999
$ws.Range("H73").Value = 146200.58
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 203880.8
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 611642.3999999999
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -613514.3999999999
$ws.Range("H129").Value = 1450.0588
$ws.Range("I129").Value = 804.25
$ws.Range("K129").Value = 2412.75
$ws.Range("M129").Value = 2587.25
$ws.Range("H132").Value = 2184.2222
$ws.Range("I132").Value = 1421.0667
$ws.Range("K132").Value = 4263.2001
$ws.Range("M132").Value = -1733.2001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3616.0571
$ws.Range("I2").Value = 3212.3333
$ws.Range("K2").Value = 3212.3333
$ws.Range("M2").Value = -3099.3333
$ws.Range("H32").Value = 2476.5217
$ws.Range("I32").Value = 2364.0447
$ws.Range("K32").Value = 2364.0447
$ws.Range("M32").Value = -2077.0447
$ws.Range("H36").Value = 14999.75
$ws.Range("I36").Value = 18999.5
$ws.Range("J36").Value = 11000
$ws.Range("K36").Value = 18999.5
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = -18653.5
$ws.Range("N36").Value = -11692
$ws.Range("H45").Value = 5077.6665
$ws.Range("J45").Value = 8049.5
$ws.Range("L45").Value = 8049.5
$ws.Range("N45").Value = -8803.5
$ws.Range("H110").Value = 2500500
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H116").Value = 3616.0571
$ws.Range("I116").Value = 3212.3333
$ws.Range("K116").Value = 3212.3333
$ws.Range("M116").Value = -918.3332999999998
$ws.Range("H122").Value = 5320.5
$ws.Range("I122").Value = 4457.857
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 13373.571
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -10923.571
$ws.Range("N122").Value = -26900.0005

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3616.0571
$ws.Range("I3").Value = 3212.3333
$ws.Range("K3").Value = 3212.3333
$ws.Range("M3").Value = -3098.3333
$ws.Range("H86").Value = 1161.3334
$ws.Range("I86").Value = 865.7222
$ws.Range("K86").Value = 865.7222
$ws.Range("M86").Value = 257.2778
$ws.Range("H89").Value = 1161.3334
$ws.Range("I89").Value = 865.7222
$ws.Range("K89").Value = 4328.611
$ws.Range("M89").Value = 1287.389

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4176.923
$ws.Range("I16").Value = 2512.625
$ws.Range("K16").Value = 2512.625
$ws.Range("M16").Value = -2225.625
$ws.Range("H58").Value = 438419.9
$ws.Range("I58").Value = 715770.6
$ws.Range("K58").Value = 715770.6
$ws.Range("M58").Value = -715567.6
$ws.Range("H105").Value = 1311.7391
$ws.Range("I105").Value = 1272.1578
$ws.Range("J105").Value = 1499.75
$ws.Range("K105").Value = 1272.1578
$ws.Range("L105").Value = 1499.75
$ws.Range("M105").Value = 474.8422
$ws.Range("N105").Value = -4993.75
$ws.Range("H113").Value = 4176.923
$ws.Range("I113").Value = 2512.625
$ws.Range("K113").Value = 2512.625
$ws.Range("M113").Value = -342.625
$ws.Range("H134").Value = 5003.5947
$ws.Range("I134").Value = 2868.2222
$ws.Range("J134").Value = 10769.1
$ws.Range("K134").Value = 8604.6666
$ws.Range("L134").Value = 32307.3
$ws.Range("M134").Value = -6069.6666
$ws.Range("N134").Value = -37377.3
$ws.Range("H136").Value = 438419.9
$ws.Range("I136").Value = 715770.6
$ws.Range("K136").Value = 2147311.8
$ws.Range("M136").Value = -2144761.8

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1173.3334
$ws.Range("J7").Value = 1500
$ws.Range("L7").Value = 4500
$ws.Range("N7").Value = -4724
$ws.Range("H39").Value = 6178.643
$ws.Range("J39").Value = 11854.857
$ws.Range("L39").Value = 35564.571
$ws.Range("N39").Value = -36152.571
$ws.Range("H41").Value = 500
$ws.Range("I41").Value = 500
$ws.Range("K41").Value = 1500
$ws.Range("H92").Value = 1025.6364
$ws.Range("I92").Value = 538.6667
$ws.Range("J92").Value = 1610
$ws.Range("K92").Value = 1616.0001
$ws.Range("L92").Value = 4830
$ws.Range("M92").Value = -368.0001
$ws.Range("N92").Value = -7326
$ws.Range("H132").Value = 4480.8
$ws.Range("I132").Value = 1702
$ws.Range("J132").Value = 6333.3335
$ws.Range("K132").Value = 15318
$ws.Range("L132").Value = 57000.0015
$ws.Range("M132").Value = -12788
$ws.Range("N132").Value = -62060.0015

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 27474.875
$ws.Range("J43").Value = 38960
$ws.Range("L43").Value = 38960
$ws.Range("N43").Value = -39262
$ws.Range("H97").Value = 9082.25
$ws.Range("I97").Value = 9082.25
$ws.Range("K97").Value = 9082.25
$ws.Range("M97").Value = -8586.25
$ws.Range("H122").Value = 11386.23
$ws.Range("I122").Value = 11003
$ws.Range("J122").Value = 11833.333
$ws.Range("K122").Value = 33009
$ws.Range("L122").Value = 35499.999
$ws.Range("M122").Value = -30559
$ws.Range("N122").Value = -40399.999
$ws.Range("H126").Value = 3871.1428
$ws.Range("I126").Value = 3433
$ws.Range("K126").Value = 10299
$ws.Range("M126").Value = -7829

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1671958
$ws.Range("I40").Value = 2503943.2
$ws.Range("K40").Value = 2503943.2
$ws.Range("M40").Value = -2503807.2
$ws.Range("H61").Value = 4806.4287
$ws.Range("I61").Value = 3529
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 3529
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -3327
$ws.Range("N61").Value = -8404
$ws.Range("H68").Value = 2190.5715
$ws.Range("I68").Value = 1766.8
$ws.Range("J68").Value = 3250
$ws.Range("K68").Value = 1766.8
$ws.Range("L68").Value = 3250
$ws.Range("M68").Value = -1017.8
$ws.Range("N68").Value = -4748
$ws.Range("H71").Value = 2190.5715
$ws.Range("I71").Value = 1766.8
$ws.Range("J71").Value = 3250
$ws.Range("K71").Value = 8834
$ws.Range("L71").Value = 16250
$ws.Range("M71").Value = -5090
$ws.Range("N71").Value = -23738
$ws.Range("H93").Value = 1028.8
$ws.Range("I93").Value = 956.38464
$ws.Range("K93").Value = 956.38464
$ws.Range("M93").Value = 291.61536
$ws.Range("H113").Value = 4806.4287
$ws.Range("I113").Value = 3529
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 3529
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -1359
$ws.Range("N113").Value = -12340
$ws.Range("H132").Value = 6151.6924
$ws.Range("I132").Value = 5295.3335
$ws.Range("J132").Value = 6885.7144
$ws.Range("K132").Value = 15886.0005
$ws.Range("L132").Value = 20657.1432
$ws.Range("M132").Value = -13356.0005
$ws.Range("N132").Value = -25717.1432

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 725.61536
$ws.Range("J2").Value = 867.6667
$ws.Range("L2").Value = 867.6667
$ws.Range("N2").Value = -1091.6667
$ws.Range("H107").Value = 37028.895
$ws.Range("J107").Value = 1070.2858
$ws.Range("L107").Value = 3210.8574
$ws.Range("N107").Value = -7050.857400000001
$ws.Range("H122").Value = 41671188
$ws.Range("I122").Value = 66669640
$ws.Range("K122").Value = 200008920
$ws.Range("M122").Value = -200006470
$ws.Range("H132").Value = 2742.0625
$ws.Range("I132").Value = 1913.3077
$ws.Range("J132").Value = 6333.3335
$ws.Range("K132").Value = 5739.9231
$ws.Range("L132").Value = 19000.0005
$ws.Range("M132").Value = -3209.9231
$ws.Range("N132").Value = -24060.0005

# ---- Structural changes ----
# ARM row 110: remove N110 (no HQ leve price/profit anymore)
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("N110").ClearContents()

# CUL row 41: add M41 = -1162 (new HQ-less profit value)
$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("M41").Value = -1162
